$wb = $excel.ActiveWorkbook

# Sheet: 展览
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F8").Value = 935  # was 932
$ws.Range("F10").Value = 95  # was 94
$ws.Range("F11").Value = 100  # was 99
$ws.Range("F12").Value = 654  # was 653
$ws.Range("F13").Value = 915  # was 913
$ws.Range("F14").Value = 1798  # was 1797
$ws.Range("F15").Value = 3866  # was 3850
$ws.Range("F16").Value = 1148  # was 1144
$ws.Range("F18").Value = 2575  # was 2567
$ws.Range("F19").Value = 677  # was 676
$ws.Range("F20").Value = 1072  # was 1070
$ws.Range("F21").Value = 3556  # was 3547
$ws.Range("F22").Value = 748  # was 744
$ws.Range("F25").Value = 2208  # was 2202
$ws.Range("F26").Value = 108  # was 107
$ws.Range("F27").Value = 830  # was 825
$ws.Range("F29").Value = 322  # was 284
$ws.Range("F30").Value = 194  # was 190
$ws.Range("F32").Value = 1324  # was 1319
$ws.Range("F33").Value = 1935  # was 1929
$ws.Range("F34").Value = 487  # was 484
$ws.Range("F35").Value = 30  # was 25
$ws.Range("F38").Value = 274  # was 269
$ws.Range("F39").Value = 11  # was 10
$ws.Range("F42").Value = 76  # was 74

# Sheet: 演出
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F11").Value = 18  # was 17

# Sheet: 本地生活
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 416  # was 414

# Sheet: 全部类型
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 416  # was 414
$ws.Range("F7").Value = 935  # was 932
$ws.Range("F10").Value = 95  # was 94
$ws.Range("F13").Value = 100  # was 99
$ws.Range("F15").Value = 915  # was 913
$ws.Range("F16").Value = 1798  # was 1797
$ws.Range("F17").Value = 3866  # was 3850
$ws.Range("F18").Value = 1148  # was 1144
$ws.Range("F21").Value = 2575  # was 2568
$ws.Range("F23").Value = 1072  # was 1070
$ws.Range("F24").Value = 3556  # was 3547
$ws.Range("F25").Value = 748  # was 744
$ws.Range("F29").Value = 2208  # was 2202
$ws.Range("F32").Value = 18  # was 17
$ws.Range("F33").Value = 108  # was 107
$ws.Range("F35").Value = 830  # was 825
$ws.Range("F37").Value = 322  # was 285
$ws.Range("F38").Value = 194  # was 190
$ws.Range("F41").Value = 1324  # was 1319
$ws.Range("F42").Value = 1935  # was 1929
$ws.Range("F44").Value = 487  # was 484
$ws.Range("F46").Value = 274  # was 270
$ws.Range("F49").Value = 76  # was 74
